$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Header (first-page header): BTec_Logo-Orange inline picture.
# InlineShape.Name writes straight through on header stories, so a
# direct rename is enough here.
$btecLogo = $sec.Headers.Item(2).Range.InlineShapes.Item(1)
$btecLogo.Name = "image1.jpg"

# --- Footers (first-page + default): Pearson Edexcel logo inline picture.
# InlineShape.Name is a no-op when the shape lives in a footer story, so
# round-trip through a floating Shape (whose .Name setter does stick) and
# convert back to an inline picture, which restores wp:inline cleanly.
$pearsonFirst = $sec.Footers.Item(2).Range.InlineShapes.Item(1)
$pearsonFirstShape = $pearsonFirst.ConvertToShape()
$pearsonFirstShape.Name = "image2.png"
$pearsonFirstShape.ConvertToInlineShape() | Out-Null

$pearsonDefault = $sec.Footers.Item(1).Range.InlineShapes.Item(1)
$pearsonDefaultShape = $pearsonDefault.ConvertToShape()
$pearsonDefaultShape.Name = "image2.png"
$pearsonDefaultShape.ConvertToInlineShape() | Out-Null
